# Fragen-Antwort Katalog: update gemaess 2. Besprechung mit Hr. Lange.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was filtered down to only "open" rows (rows 2-15 were hidden by
# the autofilter). Clear the filter criteria and unhide those rows again,
# while keeping the autofilter range itself in place.
$ws.ShowAllData()

# Copy the existing question/answer formatting from row 15 down into the
# five rows that gain new Frage/Antwort content (16-20), so the new cells
# pick up the same wrap/alignment styling used by the rest of the table.
$ws.Range("C15").Copy()
$ws.Range("C16:C20").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("D16:D20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Status column: rows 16-20 flip from "open" to "closed".
$ws.Range("B16").Value = "closed"
$ws.Range("B17").Value = "closed"
$ws.Range("B18").Value = "closed"
$ws.Range("B19").Value = "closed"
$ws.Range("B20").Value = "closed"

# Frage column (16-20) - new questions added by Hr. Lange's feedback round.
$ws.Range("C16").Value = "Zugservice/-komposition: muss eine Lokomotive in der Zugkomposition abgebildet werden?"
$ws.Range("C17").Value = "Wie soll eine Strecke am besten mit An- & Abfahrtszeiten abgebildet werden: Zeitangaben oder Dauerangaben (Dauerangaben ermöglichen eine schnellere Anpassung der Zeiten sowohl für die Erstellung der Zugservices wie auch bei Verspätungen)?"
$ws.Range("C18").Value = "Sind Vorgaben für die Platznummerierung vorhangen (à la BLS)?"
$ws.Range("C19").Value = "Soll die Applikation für z.B. 3 zu reservierende Plätze Vorschläge unterbreiten?"
$ws.Range("C20").Value = "Besteht ein Wagen aus Abteilen? Oder ist diese Differenzierung nicht notwendig?"

# Antwort column (16-20) - matching answers.
$ws.Range("D16").Value = "Nein, ist für die Sitzplatzreservation irrelevant"
$ws.Range("D17").Value = "Dauerangaben bzw. Zeitintervalle."
$ws.Range("D18").Value = "Keine vorgesehen. Admin soll selber entscheiden können."
$ws.Range("D19").Value = "Gute Idee, verbesserte UX. Option!"
$ws.Range("D20").Value = "Keine Abteile."

# Row 17 needs extra height to fit its longer question; row 20 is new and
# uses the standard question-row height.
$ws.Rows(17).RowHeight = 42
$ws.Rows(20).RowHeight = 17

# Move the active selection to reflect where editing left off.
$ws.Range("D21").Select()
